# NYPD CompStat weekly report refresh: new crime data collected.
# Updates the report header (volume/number + the covered-week date range)
# and refreshes every Week-to-Date / 28-Day / Year-to-Date / 2-Year figure
# (and their derived % Chg columns) on rows 15-29 of the CompStat_1 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  31" -> "...32", and the reporting
# week date range 7/31/2023-8/6/2023 -> 8/7/2023-8/13/2023.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# ---------------------------------------------------------------------
# A few cells flip from the "no data" text placeholder to a real number
# this week. Set their NumberFormat to match their sibling cells so the
# underlying style matches a genuine numeric cell, then assign values.
# ---------------------------------------------------------------------
$intFmt = "#,##0"
$pctFmt = '#,##0.0;"-"#,##0.0'

$ws.Range("C18").NumberFormat = $intFmt
$ws.Range("D28").NumberFormat = $intFmt
$ws.Range("G28").NumberFormat = $intFmt
$ws.Range("D29").NumberFormat = $intFmt
$ws.Range("G29").NumberFormat = $intFmt

$ws.Range("E28").NumberFormat = $pctFmt
$ws.Range("H28").NumberFormat = $pctFmt
$ws.Range("E29").NumberFormat = $pctFmt
$ws.Range("H29").NumberFormat = $pctFmt

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("L15").Value = 14.285714285714

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -36.842105263157
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 111
$ws.Range("L16").Value = -15.909090909090
$ws.Range("M16").Value = -43.076923076923
$ws.Range("N16").Value = -86.618444846292

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 18.75
$ws.Range("I17").Value = 120
$ws.Range("J17").Value = 114
$ws.Range("K17").Value = 5.263157894736
$ws.Range("L17").Value = 17.647058823529
$ws.Range("M17").Value = 90.476190476190
$ws.Range("N17").Value = -61.904761904761

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 111
$ws.Range("J18").Value = 110
$ws.Range("K18").Value = 0.909090909090
$ws.Range("L18").Value = 58.571428571428
$ws.Range("M18").Value = 38.75
$ws.Range("N18").Value = -82.870370370370

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 315
$ws.Range("J19").Value = 309
$ws.Range("K19").Value = 1.941747572815
$ws.Range("L19").Value = 22.093023255814
$ws.Range("M19").Value = -0.943396226415
$ws.Range("N19").Value = -50.78125

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 24.074074074074
$ws.Range("L20").Value = 24.074074074074
$ws.Range("M20").Value = 179.166666666667
$ws.Range("N20").Value = -89.28

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 100
$ws.Range("H21").Value = 4.166666666666
$ws.Range("I21").Value = 697
$ws.Range("J21").Value = 704
$ws.Range("K21").Value = -0.994318181818
$ws.Range("L21").Value = 20.172413793103
$ws.Range("M21").Value = 11.698717948717
$ws.Range("N21").Value = -75.274920184462

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 13
$ws.Range("K22").Value = -18.75
$ws.Range("L22").Value = -18.75
$ws.Range("M22").Value = -18.75

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 72
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = -2.702702702702
$ws.Range("M23").Value = 46.938775510204

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -18.918918918918
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = -35.947712418300
$ws.Range("I24").Value = 1038
$ws.Range("J24").Value = 1121
$ws.Range("K24").Value = -7.404103479036
$ws.Range("L24").Value = 62.441314553990
$ws.Range("M24").Value = 67.689822294022

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 80
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 186
$ws.Range("J25").Value = 190
$ws.Range("K25").Value = -2.105263157894
$ws.Range("L25").Value = 15.527950310559
$ws.Range("M25").Value = -7.462686567164

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("L26").Value = 10

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 31
$ws.Range("K27").Value = -3.125
$ws.Range("L27").Value = -6.060606060606

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 25
$ws.Range("N28").Value = -84.848484848484

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 33.333333333333
$ws.Range("N29").Value = -87.5
